$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ================================================================
# The source feed re-ordered a handful of fixtures that share the
# same kickoff date/time; re-apply rows 50/51, 89/90, 100/102/104
# and 117/118 in their corrected order, then append the three new
# fixtures (rows 140-142) added by todays update.
# ================================================================

# --- Rows 50 and 51 swap places ---
$r50_B = $ws.Range("B50").Value()
$r50_C = $ws.Range("C50").Value()
$r50_D = $ws.Range("D50").Value()
$r50_E = $ws.Range("E50").Value()
$r50_F = $ws.Range("F50").Value()
$r50_G = $ws.Range("G50").Value()
$r50_H = $ws.Range("H50").Value()
$r50_I = $ws.Range("I50").Value()
$r50_J = $ws.Range("J50").Value()
$r50_K = $ws.Range("K50").Value()
$r50_L = $ws.Range("L50").Value()
$r50_M = $ws.Range("M50").Value()
$r50_N = $ws.Range("N50").Value()
$r50_O = $ws.Range("O50").Value()
$r50_P = $ws.Range("P50").Value()
$r50_Q = $ws.Range("Q50").Value()
$r50_R = $ws.Range("R50").Value()
$r50_S = $ws.Range("S50").Value()
$r50_T = $ws.Range("T50").Value()
$r50_U = $ws.Range("U50").Value()
$r50_V = $ws.Range("V50").Value()
$r50_W = $ws.Range("W50").Value()
$r50_X = $ws.Range("X50").Value()
$r50_Y = $ws.Range("Y50").Value()
$r50_Z = $ws.Range("Z50").Value()
$r50_AA = $ws.Range("AA50").Value()
$r50_AB = $ws.Range("AB50").Value()
$r50_AC = $ws.Range("AC50").Value()
$r51_B = $ws.Range("B51").Value()
$r51_C = $ws.Range("C51").Value()
$r51_D = $ws.Range("D51").Value()
$r51_E = $ws.Range("E51").Value()
$r51_F = $ws.Range("F51").Value()
$r51_G = $ws.Range("G51").Value()
$r51_H = $ws.Range("H51").Value()
$r51_I = $ws.Range("I51").Value()
$r51_J = $ws.Range("J51").Value()
$r51_K = $ws.Range("K51").Value()
$r51_L = $ws.Range("L51").Value()
$r51_M = $ws.Range("M51").Value()
$r51_N = $ws.Range("N51").Value()
$r51_O = $ws.Range("O51").Value()
$r51_P = $ws.Range("P51").Value()
$r51_Q = $ws.Range("Q51").Value()
$r51_R = $ws.Range("R51").Value()
$r51_S = $ws.Range("S51").Value()
$r51_T = $ws.Range("T51").Value()
$r51_U = $ws.Range("U51").Value()
$r51_V = $ws.Range("V51").Value()
$r51_W = $ws.Range("W51").Value()
$r51_X = $ws.Range("X51").Value()
$r51_Y = $ws.Range("Y51").Value()
$r51_Z = $ws.Range("Z51").Value()
$r51_AA = $ws.Range("AA51").Value()
$r51_AB = $ws.Range("AB51").Value()
$r51_AC = $ws.Range("AC51").Value()
$ws.Range("B50").Value = $r51_B
$ws.Range("C50").Value = $r51_C
$ws.Range("D50").Value = $r51_D
$ws.Range("E50").Value = $r51_E
$ws.Range("F50").Value = $r51_F
$ws.Range("G50").Value = $r51_G
$ws.Range("H50").Value = $r51_H
$ws.Range("I50").Value = $r51_I
$ws.Range("J50").Value = $r51_J
$ws.Range("K50").Value = $r51_K
$ws.Range("L50").Value = $r51_L
$ws.Range("M50").Value = $r51_M
$ws.Range("N50").Value = $r51_N
$ws.Range("O50").Value = $r51_O
$ws.Range("P50").Value = $r51_P
$ws.Range("Q50").Value = $r51_Q
$ws.Range("R50").Value = $r51_R
$ws.Range("S50").Value = $r51_S
$ws.Range("T50").Value = $r51_T
$ws.Range("U50").Value = $r51_U
$ws.Range("V50").Value = $r51_V
$ws.Range("W50").Value = $r51_W
$ws.Range("X50").Value = $r51_X
$ws.Range("Y50").Value = $r51_Y
$ws.Range("Z50").Value = $r51_Z
$ws.Range("AA50").Value = $r51_AA
$ws.Range("AB50").Value = $r51_AB
$ws.Range("AC50").Value = $r51_AC
$ws.Range("B51").Value = $r50_B
$ws.Range("C51").Value = $r50_C
$ws.Range("D51").Value = $r50_D
$ws.Range("E51").Value = $r50_E
$ws.Range("F51").Value = $r50_F
$ws.Range("G51").Value = $r50_G
$ws.Range("H51").Value = $r50_H
$ws.Range("I51").Value = $r50_I
$ws.Range("J51").Value = $r50_J
$ws.Range("K51").Value = $r50_K
$ws.Range("L51").Value = $r50_L
$ws.Range("M51").Value = $r50_M
$ws.Range("N51").Value = $r50_N
$ws.Range("O51").Value = $r50_O
$ws.Range("P51").Value = $r50_P
$ws.Range("Q51").Value = $r50_Q
$ws.Range("R51").Value = $r50_R
$ws.Range("S51").Value = $r50_S
$ws.Range("T51").Value = $r50_T
$ws.Range("U51").Value = $r50_U
$ws.Range("V51").Value = $r50_V
$ws.Range("W51").Value = $r50_W
$ws.Range("X51").Value = $r50_X
$ws.Range("Y51").Value = $r50_Y
$ws.Range("Z51").Value = $r50_Z
$ws.Range("AA51").Value = $r50_AA
$ws.Range("AB51").Value = $r50_AB
$ws.Range("AC51").Value = $r50_AC

# --- Rows 89 and 90 swap places ---
$r89_B = $ws.Range("B89").Value()
$r89_C = $ws.Range("C89").Value()
$r89_D = $ws.Range("D89").Value()
$r89_E = $ws.Range("E89").Value()
$r89_F = $ws.Range("F89").Value()
$r89_G = $ws.Range("G89").Value()
$r89_H = $ws.Range("H89").Value()
$r89_I = $ws.Range("I89").Value()
$r89_J = $ws.Range("J89").Value()
$r89_K = $ws.Range("K89").Value()
$r89_L = $ws.Range("L89").Value()
$r89_M = $ws.Range("M89").Value()
$r89_N = $ws.Range("N89").Value()
$r89_O = $ws.Range("O89").Value()
$r89_P = $ws.Range("P89").Value()
$r89_Q = $ws.Range("Q89").Value()
$r89_R = $ws.Range("R89").Value()
$r89_S = $ws.Range("S89").Value()
$r89_T = $ws.Range("T89").Value()
$r89_U = $ws.Range("U89").Value()
$r89_V = $ws.Range("V89").Value()
$r89_W = $ws.Range("W89").Value()
$r89_X = $ws.Range("X89").Value()
$r89_Y = $ws.Range("Y89").Value()
$r89_Z = $ws.Range("Z89").Value()
$r89_AA = $ws.Range("AA89").Value()
$r89_AB = $ws.Range("AB89").Value()
$r89_AC = $ws.Range("AC89").Value()
$r90_B = $ws.Range("B90").Value()
$r90_C = $ws.Range("C90").Value()
$r90_D = $ws.Range("D90").Value()
$r90_E = $ws.Range("E90").Value()
$r90_F = $ws.Range("F90").Value()
$r90_G = $ws.Range("G90").Value()
$r90_H = $ws.Range("H90").Value()
$r90_I = $ws.Range("I90").Value()
$r90_J = $ws.Range("J90").Value()
$r90_K = $ws.Range("K90").Value()
$r90_L = $ws.Range("L90").Value()
$r90_M = $ws.Range("M90").Value()
$r90_N = $ws.Range("N90").Value()
$r90_O = $ws.Range("O90").Value()
$r90_P = $ws.Range("P90").Value()
$r90_Q = $ws.Range("Q90").Value()
$r90_R = $ws.Range("R90").Value()
$r90_S = $ws.Range("S90").Value()
$r90_T = $ws.Range("T90").Value()
$r90_U = $ws.Range("U90").Value()
$r90_V = $ws.Range("V90").Value()
$r90_W = $ws.Range("W90").Value()
$r90_X = $ws.Range("X90").Value()
$r90_Y = $ws.Range("Y90").Value()
$r90_Z = $ws.Range("Z90").Value()
$r90_AA = $ws.Range("AA90").Value()
$r90_AB = $ws.Range("AB90").Value()
$r90_AC = $ws.Range("AC90").Value()
$ws.Range("B89").Value = $r90_B
$ws.Range("C89").Value = $r90_C
$ws.Range("D89").Value = $r90_D
$ws.Range("E89").Value = $r90_E
$ws.Range("F89").Value = $r90_F
$ws.Range("G89").Value = $r90_G
$ws.Range("H89").Value = $r90_H
$ws.Range("I89").Value = $r90_I
$ws.Range("J89").Value = $r90_J
$ws.Range("K89").Value = $r90_K
$ws.Range("L89").Value = $r90_L
$ws.Range("M89").Value = $r90_M
$ws.Range("N89").Value = $r90_N
$ws.Range("O89").Value = $r90_O
$ws.Range("P89").Value = $r90_P
$ws.Range("Q89").Value = $r90_Q
$ws.Range("R89").Value = $r90_R
$ws.Range("S89").Value = $r90_S
$ws.Range("T89").Value = $r90_T
$ws.Range("U89").Value = $r90_U
$ws.Range("V89").Value = $r90_V
$ws.Range("W89").Value = $r90_W
$ws.Range("X89").Value = $r90_X
$ws.Range("Y89").Value = $r90_Y
$ws.Range("Z89").Value = $r90_Z
$ws.Range("AA89").Value = $r90_AA
$ws.Range("AB89").Value = $r90_AB
$ws.Range("AC89").Value = $r90_AC
$ws.Range("B90").Value = $r89_B
$ws.Range("C90").Value = $r89_C
$ws.Range("D90").Value = $r89_D
$ws.Range("E90").Value = $r89_E
$ws.Range("F90").Value = $r89_F
$ws.Range("G90").Value = $r89_G
$ws.Range("H90").Value = $r89_H
$ws.Range("I90").Value = $r89_I
$ws.Range("J90").Value = $r89_J
$ws.Range("K90").Value = $r89_K
$ws.Range("L90").Value = $r89_L
$ws.Range("M90").Value = $r89_M
$ws.Range("N90").Value = $r89_N
$ws.Range("O90").Value = $r89_O
$ws.Range("P90").Value = $r89_P
$ws.Range("Q90").Value = $r89_Q
$ws.Range("R90").Value = $r89_R
$ws.Range("S90").Value = $r89_S
$ws.Range("T90").Value = $r89_T
$ws.Range("U90").Value = $r89_U
$ws.Range("V90").Value = $r89_V
$ws.Range("W90").Value = $r89_W
$ws.Range("X90").Value = $r89_X
$ws.Range("Y90").Value = $r89_Y
$ws.Range("Z90").Value = $r89_Z
$ws.Range("AA90").Value = $r89_AA
$ws.Range("AB90").Value = $r89_AB
$ws.Range("AC90").Value = $r89_AC

# --- Rows 100, 102, 104 rotate: 100<-102, 102<-104, 104<-100 ---
$r100_B = $ws.Range("B100").Value()
$r100_C = $ws.Range("C100").Value()
$r100_D = $ws.Range("D100").Value()
$r100_E = $ws.Range("E100").Value()
$r100_F = $ws.Range("F100").Value()
$r100_G = $ws.Range("G100").Value()
$r100_H = $ws.Range("H100").Value()
$r100_I = $ws.Range("I100").Value()
$r100_J = $ws.Range("J100").Value()
$r100_K = $ws.Range("K100").Value()
$r100_L = $ws.Range("L100").Value()
$r100_M = $ws.Range("M100").Value()
$r100_N = $ws.Range("N100").Value()
$r100_O = $ws.Range("O100").Value()
$r100_P = $ws.Range("P100").Value()
$r100_Q = $ws.Range("Q100").Value()
$r100_R = $ws.Range("R100").Value()
$r100_S = $ws.Range("S100").Value()
$r100_T = $ws.Range("T100").Value()
$r100_U = $ws.Range("U100").Value()
$r100_V = $ws.Range("V100").Value()
$r100_W = $ws.Range("W100").Value()
$r100_X = $ws.Range("X100").Value()
$r100_Y = $ws.Range("Y100").Value()
$r100_Z = $ws.Range("Z100").Value()
$r100_AA = $ws.Range("AA100").Value()
$r100_AB = $ws.Range("AB100").Value()
$r100_AC = $ws.Range("AC100").Value()
$r102_B = $ws.Range("B102").Value()
$r102_C = $ws.Range("C102").Value()
$r102_D = $ws.Range("D102").Value()
$r102_E = $ws.Range("E102").Value()
$r102_F = $ws.Range("F102").Value()
$r102_G = $ws.Range("G102").Value()
$r102_H = $ws.Range("H102").Value()
$r102_I = $ws.Range("I102").Value()
$r102_J = $ws.Range("J102").Value()
$r102_K = $ws.Range("K102").Value()
$r102_L = $ws.Range("L102").Value()
$r102_M = $ws.Range("M102").Value()
$r102_N = $ws.Range("N102").Value()
$r102_O = $ws.Range("O102").Value()
$r102_P = $ws.Range("P102").Value()
$r102_Q = $ws.Range("Q102").Value()
$r102_R = $ws.Range("R102").Value()
$r102_S = $ws.Range("S102").Value()
$r102_T = $ws.Range("T102").Value()
$r102_U = $ws.Range("U102").Value()
$r102_V = $ws.Range("V102").Value()
$r102_W = $ws.Range("W102").Value()
$r102_X = $ws.Range("X102").Value()
$r102_Y = $ws.Range("Y102").Value()
$r102_Z = $ws.Range("Z102").Value()
$r102_AA = $ws.Range("AA102").Value()
$r102_AB = $ws.Range("AB102").Value()
$r102_AC = $ws.Range("AC102").Value()
$r104_B = $ws.Range("B104").Value()
$r104_C = $ws.Range("C104").Value()
$r104_D = $ws.Range("D104").Value()
$r104_E = $ws.Range("E104").Value()
$r104_F = $ws.Range("F104").Value()
$r104_G = $ws.Range("G104").Value()
$r104_H = $ws.Range("H104").Value()
$r104_I = $ws.Range("I104").Value()
$r104_J = $ws.Range("J104").Value()
$r104_K = $ws.Range("K104").Value()
$r104_L = $ws.Range("L104").Value()
$r104_M = $ws.Range("M104").Value()
$r104_N = $ws.Range("N104").Value()
$r104_O = $ws.Range("O104").Value()
$r104_P = $ws.Range("P104").Value()
$r104_Q = $ws.Range("Q104").Value()
$r104_R = $ws.Range("R104").Value()
$r104_S = $ws.Range("S104").Value()
$r104_T = $ws.Range("T104").Value()
$r104_U = $ws.Range("U104").Value()
$r104_V = $ws.Range("V104").Value()
$r104_W = $ws.Range("W104").Value()
$r104_X = $ws.Range("X104").Value()
$r104_Y = $ws.Range("Y104").Value()
$r104_Z = $ws.Range("Z104").Value()
$r104_AA = $ws.Range("AA104").Value()
$r104_AB = $ws.Range("AB104").Value()
$r104_AC = $ws.Range("AC104").Value()
$ws.Range("B100").Value = $r102_B
$ws.Range("C100").Value = $r102_C
$ws.Range("D100").Value = $r102_D
$ws.Range("E100").Value = $r102_E
$ws.Range("F100").Value = $r102_F
$ws.Range("G100").Value = $r102_G
$ws.Range("H100").Value = $r102_H
$ws.Range("I100").Value = $r102_I
$ws.Range("J100").Value = $r102_J
$ws.Range("K100").Value = $r102_K
$ws.Range("L100").Value = $r102_L
$ws.Range("M100").Value = $r102_M
$ws.Range("N100").Value = $r102_N
$ws.Range("O100").Value = $r102_O
$ws.Range("P100").Value = $r102_P
$ws.Range("Q100").Value = $r102_Q
$ws.Range("R100").Value = $r102_R
$ws.Range("S100").Value = $r102_S
$ws.Range("T100").Value = $r102_T
$ws.Range("U100").Value = $r102_U
$ws.Range("V100").Value = $r102_V
$ws.Range("W100").Value = $r102_W
$ws.Range("X100").Value = $r102_X
$ws.Range("Y100").Value = $r102_Y
$ws.Range("Z100").Value = $r102_Z
$ws.Range("AA100").Value = $r102_AA
$ws.Range("AB100").Value = $r102_AB
$ws.Range("AC100").Value = $r102_AC
$ws.Range("B102").Value = $r104_B
$ws.Range("C102").Value = $r104_C
$ws.Range("D102").Value = $r104_D
$ws.Range("E102").Value = $r104_E
$ws.Range("F102").Value = $r104_F
$ws.Range("G102").Value = $r104_G
$ws.Range("H102").Value = $r104_H
$ws.Range("I102").Value = $r104_I
$ws.Range("J102").Value = $r104_J
$ws.Range("K102").Value = $r104_K
$ws.Range("L102").Value = $r104_L
$ws.Range("M102").Value = $r104_M
$ws.Range("N102").Value = $r104_N
$ws.Range("O102").Value = $r104_O
$ws.Range("P102").Value = $r104_P
$ws.Range("Q102").Value = $r104_Q
$ws.Range("R102").Value = $r104_R
$ws.Range("S102").Value = $r104_S
$ws.Range("T102").Value = $r104_T
$ws.Range("U102").Value = $r104_U
$ws.Range("V102").Value = $r104_V
$ws.Range("W102").Value = $r104_W
$ws.Range("X102").Value = $r104_X
$ws.Range("Y102").Value = $r104_Y
$ws.Range("Z102").Value = $r104_Z
$ws.Range("AA102").Value = $r104_AA
$ws.Range("AB102").Value = $r104_AB
$ws.Range("AC102").Value = $r104_AC
$ws.Range("B104").Value = $r100_B
$ws.Range("C104").Value = $r100_C
$ws.Range("D104").Value = $r100_D
$ws.Range("E104").Value = $r100_E
$ws.Range("F104").Value = $r100_F
$ws.Range("G104").Value = $r100_G
$ws.Range("H104").Value = $r100_H
$ws.Range("I104").Value = $r100_I
$ws.Range("J104").Value = $r100_J
$ws.Range("K104").Value = $r100_K
$ws.Range("L104").Value = $r100_L
$ws.Range("M104").Value = $r100_M
$ws.Range("N104").Value = $r100_N
$ws.Range("O104").Value = $r100_O
$ws.Range("P104").Value = $r100_P
$ws.Range("Q104").Value = $r100_Q
$ws.Range("R104").Value = $r100_R
$ws.Range("S104").Value = $r100_S
$ws.Range("T104").Value = $r100_T
$ws.Range("U104").Value = $r100_U
$ws.Range("V104").Value = $r100_V
$ws.Range("W104").Value = $r100_W
$ws.Range("X104").Value = $r100_X
$ws.Range("Y104").Value = $r100_Y
$ws.Range("Z104").Value = $r100_Z
$ws.Range("AA104").Value = $r100_AA
$ws.Range("AB104").Value = $r100_AB
$ws.Range("AC104").Value = $r100_AC

# --- Rows 117 and 118 swap places ---
$r117_B = $ws.Range("B117").Value()
$r117_C = $ws.Range("C117").Value()
$r117_D = $ws.Range("D117").Value()
$r117_E = $ws.Range("E117").Value()
$r117_F = $ws.Range("F117").Value()
$r117_G = $ws.Range("G117").Value()
$r117_H = $ws.Range("H117").Value()
$r117_I = $ws.Range("I117").Value()
$r117_J = $ws.Range("J117").Value()
$r117_K = $ws.Range("K117").Value()
$r117_L = $ws.Range("L117").Value()
$r117_M = $ws.Range("M117").Value()
$r117_N = $ws.Range("N117").Value()
$r117_O = $ws.Range("O117").Value()
$r117_P = $ws.Range("P117").Value()
$r117_Q = $ws.Range("Q117").Value()
$r117_R = $ws.Range("R117").Value()
$r117_S = $ws.Range("S117").Value()
$r117_T = $ws.Range("T117").Value()
$r117_U = $ws.Range("U117").Value()
$r117_V = $ws.Range("V117").Value()
$r117_W = $ws.Range("W117").Value()
$r117_X = $ws.Range("X117").Value()
$r117_Y = $ws.Range("Y117").Value()
$r117_Z = $ws.Range("Z117").Value()
$r117_AA = $ws.Range("AA117").Value()
$r117_AB = $ws.Range("AB117").Value()
$r117_AC = $ws.Range("AC117").Value()
$r118_B = $ws.Range("B118").Value()
$r118_C = $ws.Range("C118").Value()
$r118_D = $ws.Range("D118").Value()
$r118_E = $ws.Range("E118").Value()
$r118_F = $ws.Range("F118").Value()
$r118_G = $ws.Range("G118").Value()
$r118_H = $ws.Range("H118").Value()
$r118_I = $ws.Range("I118").Value()
$r118_J = $ws.Range("J118").Value()
$r118_K = $ws.Range("K118").Value()
$r118_L = $ws.Range("L118").Value()
$r118_M = $ws.Range("M118").Value()
$r118_N = $ws.Range("N118").Value()
$r118_O = $ws.Range("O118").Value()
$r118_P = $ws.Range("P118").Value()
$r118_Q = $ws.Range("Q118").Value()
$r118_R = $ws.Range("R118").Value()
$r118_S = $ws.Range("S118").Value()
$r118_T = $ws.Range("T118").Value()
$r118_U = $ws.Range("U118").Value()
$r118_V = $ws.Range("V118").Value()
$r118_W = $ws.Range("W118").Value()
$r118_X = $ws.Range("X118").Value()
$r118_Y = $ws.Range("Y118").Value()
$r118_Z = $ws.Range("Z118").Value()
$r118_AA = $ws.Range("AA118").Value()
$r118_AB = $ws.Range("AB118").Value()
$r118_AC = $ws.Range("AC118").Value()
$ws.Range("B117").Value = $r118_B
$ws.Range("C117").Value = $r118_C
$ws.Range("D117").Value = $r118_D
$ws.Range("E117").Value = $r118_E
$ws.Range("F117").Value = $r118_F
$ws.Range("G117").Value = $r118_G
$ws.Range("H117").Value = $r118_H
$ws.Range("I117").Value = $r118_I
$ws.Range("J117").Value = $r118_J
$ws.Range("K117").Value = $r118_K
$ws.Range("L117").Value = $r118_L
$ws.Range("M117").Value = $r118_M
$ws.Range("N117").Value = $r118_N
$ws.Range("O117").Value = $r118_O
$ws.Range("P117").Value = $r118_P
$ws.Range("Q117").Value = $r118_Q
$ws.Range("R117").Value = $r118_R
$ws.Range("S117").Value = $r118_S
$ws.Range("T117").Value = $r118_T
$ws.Range("U117").Value = $r118_U
$ws.Range("V117").Value = $r118_V
$ws.Range("W117").Value = $r118_W
$ws.Range("X117").Value = $r118_X
$ws.Range("Y117").Value = $r118_Y
$ws.Range("Z117").Value = $r118_Z
$ws.Range("AA117").Value = $r118_AA
$ws.Range("AB117").Value = $r118_AB
$ws.Range("AC117").Value = $r118_AC
$ws.Range("B118").Value = $r117_B
$ws.Range("C118").Value = $r117_C
$ws.Range("D118").Value = $r117_D
$ws.Range("E118").Value = $r117_E
$ws.Range("F118").Value = $r117_F
$ws.Range("G118").Value = $r117_G
$ws.Range("H118").Value = $r117_H
$ws.Range("I118").Value = $r117_I
$ws.Range("J118").Value = $r117_J
$ws.Range("K118").Value = $r117_K
$ws.Range("L118").Value = $r117_L
$ws.Range("M118").Value = $r117_M
$ws.Range("N118").Value = $r117_N
$ws.Range("O118").Value = $r117_O
$ws.Range("P118").Value = $r117_P
$ws.Range("Q118").Value = $r117_Q
$ws.Range("R118").Value = $r117_R
$ws.Range("S118").Value = $r117_S
$ws.Range("T118").Value = $r117_T
$ws.Range("U118").Value = $r117_U
$ws.Range("V118").Value = $r117_V
$ws.Range("W118").Value = $r117_W
$ws.Range("X118").Value = $r117_X
$ws.Range("Y118").Value = $r117_Y
$ws.Range("Z118").Value = $r117_Z
$ws.Range("AA118").Value = $r117_AA
$ws.Range("AB118").Value = $r117_AB
$ws.Range("AC118").Value = $r117_AC

# --- Append three newly played/scheduled fixtures ---
# Row 140
$ws.Range("A139").Copy()
$ws.Range("A140").PasteSpecial(-4122)
$ws.Range("E139").Copy()
$ws.Range("E140").PasteSpecial(-4122)
$ws.Range("A140").Value = 138
$ws.Range("B140").Value = 7862925
$ws.Range("C140").Value = "Lithuania A Lyga"
$ws.Range("D140").Value = "Lithuania A Lyga"
$ws.Range("E140").Value = 45396.29166666666
$ws.Range("F140").Value = "FK Dziugas Telsiai"
$ws.Range("G140").Value = "FK Siauliai"
$ws.Range("K140").Value = 3.75
$ws.Range("L140").Value = 3.1
$ws.Range("M140").Value = 1.909
$ws.Range("N140").Value = 3.75
$ws.Range("O140").Value = 3.1
$ws.Range("P140").Value = 1.909
$ws.Range("Q140").Value = 0.25
$ws.Range("R140").Value = 2.1
$ws.Range("S140").Value = 1.7
$ws.Range("T140").Value = 2.25
$ws.Range("U140").Value = 1.975
$ws.Range("V140").Value = 1.825
$ws.Range("W140").Value = 0
$ws.Range("X140").Value = 0
$ws.Range("Y140").Value = 0
$ws.Range("Z140").Value = 0
$ws.Range("AA140").Value = 0

# Row 141
$ws.Range("A139").Copy()
$ws.Range("A141").PasteSpecial(-4122)
$ws.Range("E139").Copy()
$ws.Range("E141").PasteSpecial(-4122)
$ws.Range("A141").Value = 139
$ws.Range("B141").Value = 7862046
$ws.Range("C141").Value = "Lithuania A Lyga"
$ws.Range("D141").Value = "Lithuania A Lyga"
$ws.Range("E141").Value = 45396.375
$ws.Range("F141").Value = "Panevezys"
$ws.Range("G141").Value = "FK Zalgiris Vilnius"
$ws.Range("K141").Value = 3.75
$ws.Range("L141").Value = 3.25
$ws.Range("M141").Value = 1.833
$ws.Range("N141").Value = 3.25
$ws.Range("O141").Value = 3
$ws.Range("P141").Value = 2.1
$ws.Range("Q141").Value = 0.25
$ws.Range("R141").Value = 1.95
$ws.Range("S141").Value = 1.85
$ws.Range("T141").Value = 2.25
$ws.Range("U141").Value = 2.025
$ws.Range("V141").Value = 1.775
$ws.Range("W141").Value = 0
$ws.Range("X141").Value = 0
$ws.Range("Y141").Value = 0
$ws.Range("Z141").Value = 0
$ws.Range("AA141").Value = 0

# Row 142
$ws.Range("A139").Copy()
$ws.Range("A142").PasteSpecial(-4122)
$ws.Range("E139").Copy()
$ws.Range("E142").PasteSpecial(-4122)
$ws.Range("A142").Value = 140
$ws.Range("B142").Value = 7862926
$ws.Range("C142").Value = "Lithuania A Lyga"
$ws.Range("D142").Value = "Lithuania A Lyga"
$ws.Range("E142").Value = 45396.52083333334
$ws.Range("F142").Value = "FK Dainava Alytus"
$ws.Range("G142").Value = "FK Kauno Zalgiris"
$ws.Range("K142").Value = 4.2
$ws.Range("L142").Value = 3.3
$ws.Range("M142").Value = 1.75
$ws.Range("N142").Value = 3.8
$ws.Range("O142").Value = 3.25
$ws.Range("P142").Value = 1.833
$ws.Range("Q142").Value = 0.5
$ws.Range("R142").Value = 1.925
$ws.Range("S142").Value = 1.875
$ws.Range("T142").Value = 2.25
$ws.Range("U142").Value = 1.9
$ws.Range("V142").Value = 1.9
$ws.Range("W142").Value = 0
$ws.Range("X142").Value = 0
$ws.Range("Y142").Value = 0
$ws.Range("Z142").Value = 0
$ws.Range("AA142").Value = 0

$excel.CutCopyMode = 0
